$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.520808
$ws.Range("H2").Value = 7.562424
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 10.000565
$ws.Range("N2").Value = 30.001695
$ws.Range("O2").Value = 0.6316353758144477
$ws.Range("P2").Value = 0.6316353758144477
$ws.Range("Q2").Value = 25.20950425652
$ws.Range("R2").Value = 226.88553830868
$ws.Range("S2").Value = 0.6316353758144477
$ws.Range("T2").Value = 0.6316353758144477

# Row 3 (Target cluster: FAPs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.520808
$ws.Range("H3").Value = 7.562424
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.264793333333333
$ws.Range("N3").Value = 12.79438
$ws.Range("O3").Value = 0.2693642149089528
$ws.Range("P3").Value = 0.2693642149089528
$ws.Range("Q3").Value = 10.75072515301333
$ws.Range("R3").Value = 96.75652637712
$ws.Range("S3").Value = 0.2693642149089528
$ws.Range("T3").Value = 0.2693642149089528

# Row 4 (Target cluster: sCs)
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.520808
$ws.Range("H4").Value = 7.562424
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.567455
$ws.Range("N4").Value = 4.702364999999999
$ws.Range("O4").Value = 0.09900040927659938
$ws.Range("P4").Value = 0.09900040927659938
$ws.Range("Q4").Value = 3.95125310364
$ws.Range("R4").Value = 35.56127793276
$ws.Range("S4").Value = 0.09900040927659938
$ws.Range("T4").Value = 0.09900040927659938
